# Apply "VOC 2 class Train results" edit to the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Carry over the formatting used on row 2's "Model Iteration" cell (F,
# centered+bordered) and its three computed-percentage cells (J:L,
# percent+centered+bordered) onto the two new rows before filling in
# values. The other numeric cells (G:I) keep the plain bordered style
# that the blank template rows already had.
$ws.Range("F2").Copy()
$ws.Range("F3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("J2:L2").Copy()
$ws.Range("J3:L3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J4:L4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 3: Faster RCNN / Chandra / VOC 2C vs VOC 2C ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Chandra"
$ws.Range("C3").Value = "Faster RCNN"
$ws.Range("D3").Value = "VOC 2C"
$ws.Range("E3").Value = "VOC 2C"
$ws.Range("F3").Value = 100000
$ws.Range("G3").Value = 4368
$ws.Range("H3").Value = 3264
$ws.Range("I3").Value = 24002
$ws.Range("J3").Formula = "=G3/(G3+H3)"
$ws.Range("K3").Formula = "=H3/(H3+I3)"
$ws.Range("L3").Formula = "=2*(J3*K3)/(J3+K3)"

# --- Row 4: Faster RCNN / Chandra / VOC 2C vs Mini Drone ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Chandra"
$ws.Range("C4").Value = "Faster RCNN"
$ws.Range("D4").Value = "VOC 2C"
$ws.Range("E4").Value = "Mini Drone"
$ws.Range("F4").Value = 100000
$ws.Range("G4").Value = 12120
$ws.Range("H4").Value = 10769
$ws.Range("I4").Value = 31589
$ws.Range("J4").Formula = "=G4/(G4+H4)"
$ws.Range("K4").Formula = "=H4/(H4+I4)"
$ws.Range("L4").Formula = "=2*(J4*K4)/(J4+K4)"

# Match selection shown in the saved file (cell J5 selected).
$ws.Range("J5").Select()
